$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "GIS & Geospatial Analysis Consulting",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "GIS & Geospatial Analysis Consulting`r• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels`r• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide`r• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis",
    2
)
